# Insert a new data row at row 317 (pushing the existing rows 317:438 down to
# 318:439) and populate it with a new "Ajo" price record for
# Feria Lagunitas de Puerto Montt.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 317; everything below shifts down
# by one (old row 317 -> new row 318, ..., old row 438 -> new row 439).
$ws.Rows("317:317").Insert()

# Fill in the new row 317 with the new record's data.
$ws.Range("A317").Value = 4
$ws.Range("B317").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C317").Value = "Los Lagos"
$ws.Range("D317").Value = 45009
$ws.Range("E317").Value = 10
$ws.Range("F317").Value = 100112003
$ws.Range("G317").Value = "Ajo"
$ws.Range("H317").Value = "Chino"
$ws.Range("I317").Value = "Primera"
$ws.Range("J317").Value = 240
$ws.Range("K317").Value = 20000
$ws.Range("L317").Value = 21000
$ws.Range("M317").Value = 20500
$ws.Range("N317").Value = "`$/caja 10 kilos"
$ws.Range("O317").Value = "China"
$ws.Range("P317").Value = 2050
$ws.Range("Q317").Value = 10
$ws.Range("R317").Value = "Hortaliza"

# Make sure the date cell keeps the same date/time number format used by the
# rest of column D.
$ws.Range("D317").NumberFormat = $ws.Range("D318").NumberFormat
